$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header label for CRLB column to include units
$ws.Range("F1").Value = "CRLB [Hz^2]"

# Updated "Mean estimated", "Mean error", "Variance" and "CRLB" values
# from corrected variance estimation / CRLB computation
$values = @{
    "C2" = 99707.03125
    "D2" = 292.96875
    "E2" = 95367.431640625
    "F2" = 22514.95401295208
    "F3" = 2251.495401295208
    "F4" = 225.1495401295208
    "F5" = 22.51495401295207
    "F6" = 2.251495401295208
    "F7" = 0.2251495401295208
    "F8" = 0.02251495401295208
    "F9" = 0.002251495401295208
    "C10" = 99926.7578125
    "D10" = 73.2421875
    "E10" = 27153.22706434462
    "F10" = 22514.95401295208
    "C11" = 100048.828125
    "D11" = -48.828125
    "E11" = 10596.38129340278
    "F11" = 2251.495401295208
    "C12" = 100073.2421875
    "D12" = -73.2421875
    "E12" = 5960.464477539062
    "F12" = 225.1495401295208
    "F13" = 22.51495401295207
    "F14" = 2.251495401295208
    "F15" = 0.2251495401295208
    "F16" = 0.02251495401295208
    "F17" = 0.002251495401295208
    "C18" = 99981.689453125
    "D18" = 18.310546875
    "E18" = 12790.16335805257
    "F18" = 22514.95401295208
    "C19" = 100000
    "D19" = 0
    "E19" = 1821.253034803603
    "F19" = 2251.495401295208
    "E20" = 869.2344029744467
    "F20" = 225.1495401295208
    "C21" = 99975.5859375
    "D21" = 24.4140625
    "E21" = 0
    "F21" = 22.51495401295207
    "F22" = 2.251495401295208
    "F23" = 0.2251495401295208
    "F24" = 0.02251495401295208
    "F25" = 0.002251495401295208
    "C26" = 99975.5859375
    "D26" = 24.4140625
    "E26" = 16142.92462666829
    "F26" = 22514.95401295208
    "C27" = 99986.26708984375
    "D27" = 13.73291015625
    "E27" = 1606.5314412117
    "F27" = 2251.495401295208
    "C28" = 100000
    "D28" = 0
    "E28" = 62.0881716410319
    "F28" = 225.1495401295208
    "C29" = 100001.5258789062
    "D29" = -1.52587890625
    "E29" = 54.32715018590292
    "F29" = 22.51495401295207
    "C30" = 100004.5776367188
    "D30" = -4.57763671875
    "E30" = 23.28306436538696
    "F30" = 2.251495401295208
    "F31" = 0.2251495401295208
    "F32" = 0.02251495401295208
    "F33" = 0.002251495401295208
    "C34" = 99993.51501464844
    "D34" = 6.4849853515625
    "E34" = 6645.536308901177
    "F34" = 22514.95401295208
    "C35" = 99986.26708984375
    "D35" = 13.73291015625
    "E35" = 785.1566705438826
    "F35" = 2251.495401295208
    "C36" = 99998.09265136719
    "D36" = 1.9073486328125
    "E36" = 127.5717901686827
    "F36" = 225.1495401295208
    "C37" = 100001.1444091797
    "D37" = -1.1444091796875
    "E37" = 9.862964765893089
    "F37" = 22.51495401295207
    "C38" = 100000.3814697266
    "D38" = -0.3814697265625
    "E38" = 4.042198674546348
    "F38" = 2.251495401295208
    "C39" = 99998.85559082031
    "D39" = 1.1444091796875
    "E39" = 1.455191522836685
    "F39" = 0.2251495401295208
    "F40" = 0.02251495401295208
    "F41" = 0.002251495401295208
    "C42" = 99978.82843017578
    "D42" = 21.17156982421875
    "E42" = 14842.99395492093
    "F42" = 22514.95401295208
    "C43" = 99973.0110168457
    "D43" = 26.98898315429688
    "E43" = 988.3276814232684
    "F43" = 2251.495401295208
    "C44" = 99995.04089355469
    "D44" = 4.9591064453125
    "E44" = 174.2591848596931
    "F44" = 225.1495401295208
    "E45" = 16.82565198279917
    "F45" = 22.51495401295207
    "C46" = 100000.2861022949
    "D46" = -0.286102294921875
    "E46" = 1.303609072541197
    "F46" = 2.251495401295208
    "C47" = 99999.80926513672
    "D47" = 0.19073486328125
    "E47" = 0.2425319204727809
    "F47" = 0.2251495401295208
    "C48" = 100000.1907348633
    "D48" = -0.19073486328125
    "E48" = 0.1616879469818539
    "F48" = 0.02251495401295208
    "C49" = 100000.3814697266
    "D49" = -0.3814697265625
    "E49" = 0
    "F49" = 0.002251495401295208
}

foreach ($cell in $values.Keys) {
    $ws.Range($cell).Value = $values[$cell]
}

